# Expense_Tracker.xlsx edit script
# - Swaps / rotates several same-date expense rows (B=Category, C=Amount)
# - Inserts two new expense rows ("test2" on 2025/04/28 and "test" on 2025/04/29),
#   shifting the previous tail rows (199-201) down to 200-203 as part of the
#   date-sorted insert, per the app's "save_expense_entry" auto-sort behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while guaranteeing it is stored as TEXT
# (prevents Excel's autoconvert of date-shaped strings like "2025/04/28"
# into a numeric date serial, which the source workbook never uses -
# every Date cell in this sheet is plain text).
function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---- Row 63 / 64 swap ----
$ws.Range("B63").Value = "Laundry"
$ws.Range("C63").Value = 75.81999999999999
$ws.Range("B64").Value = "Phone"
$ws.Range("C64").Value = 191.02

# ---- Row 88 / 89 swap (amount only) ----
$ws.Range("C88").Value = 119.13
$ws.Range("C89").Value = 17.42

# ---- Row 100 / 101 / 102 rotate ----
$ws.Range("B100").Value = "Toters"
$ws.Range("C100").Value = 52.66
$ws.Range("B101").Value = "Snacks"
$ws.Range("C101").Value = 12.83
$ws.Range("B102").Value = "Groceries"
$ws.Range("C102").Value = 80.29000000000001

# ---- Row 113 / 114 swap (amount only) ----
$ws.Range("C113").Value = 126.16
$ws.Range("C114").Value = 74.06999999999999

# ---- Row 126 / 127 swap ----
$ws.Range("B126").Value = "Entertainment"
$ws.Range("C126").Value = 127.43
$ws.Range("B127").Value = "Toters"
$ws.Range("C127").Value = 158.48

# ---- Row 137 / 139 swap (row 138 untouched) ----
$ws.Range("B137").Value = "Barber"
$ws.Range("C137").Value = 21.48
$ws.Range("B139").Value = "Phone"
$ws.Range("C139").Value = 146.69

# ---- Row 149 / 150 / 151 / 152 rotate ----
$ws.Range("B149").Value = "Entertainment"
$ws.Range("C149").Value = 101.09
$ws.Range("B150").Value = "Snacks"
$ws.Range("C150").Value = 6.21
$ws.Range("B151").Value = "Barber"
$ws.Range("C151").Value = 199
$ws.Range("B152").Value = "Groceries"
$ws.Range("C152").Value = 125.39

# ---- Row 162 / 163 / 164 / 165 / 166 rotate ----
$ws.Range("B162").Value = "Snacks"
$ws.Range("C162").Value = 177.11
$ws.Range("B163").Value = "Barber"
$ws.Range("C163").Value = 180.03
$ws.Range("B164").Value = "Groceries"
$ws.Range("C164").Value = 115.05
$ws.Range("B165").Value = "Restaurant"
$ws.Range("C165").Value = 105.4
$ws.Range("B166").Value = "Phone"
$ws.Range("C166").Value = 159.19

# ---- Row 173 / 174 swap (amount only) ----
$ws.Range("C173").Value = 47.34
$ws.Range("C174").Value = 111.86

# ---- Row 175 / 176 swap ----
$ws.Range("B175").Value = "Groceries"
$ws.Range("C175").Value = 72.81
$ws.Range("B176").Value = "Toters"
$ws.Range("C176").Value = 196.23

# ---- Tail restructure: insert the two new rows (sorted by date) ----
# New row 199: 2025/04/28 Restaurant 20 "test2"
Set-TextValue $ws.Range("A199") "2025/04/28"
$ws.Range("B199").Value = "Restaurant"
$ws.Range("C199").Value = 20
$ws.Range("D199").Value = "test2"

# Row 200 (was 199): 2025/04/29 Barber 69.85
# (D200 intentionally left alone: row 199 already had a blank Notes cell,
# and writing "" to it would delete the cell outright instead of keeping it blank)
Set-TextValue $ws.Range("A200") "2025/04/29"
$ws.Range("B200").Value = "Barber"
$ws.Range("C200").Value = 69.84999999999999

# New row 201: 2025/04/29 Restaurant 20 "test"
Set-TextValue $ws.Range("A201") "2025/04/29"
$ws.Range("B201").Value = "Restaurant"
$ws.Range("C201").Value = 20
$ws.Range("D201").Value = "test"

# Row 202 (was 200): 2025/04/30 Groceries 94.58
Set-TextValue $ws.Range("A202") "2025/04/30"
$ws.Range("B202").Value = "Groceries"
$ws.Range("C202").Value = 94.58

# Row 203 (was 201): 2025/04/30 Shopping 6.51
Set-TextValue $ws.Range("A203") "2025/04/30"
$ws.Range("B203").Value = "Shopping"
$ws.Range("C203").Value = 6.51

Write-Output "Edit complete."
